$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'42.320.61"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -2.75%  "
$ws.Range("D3").Value = "'2.218.04"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.93%  "
$ws.Range("D4").Value = "'1.01"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.27%  "
$ws.Range("D5").Value = "'107.66"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -9.14%  "
$ws.Range("D6").Value = "'296.74"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +12.26%  "
$ws.Range("D7").Value = "'0.620"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -3.41%  "
$ws.Range("E8").Value = "  -0.26%  "
$ws.Range("E9").Value = "  -3.05%  "
$ws.Range("E10").Value = "  -7.93%  "
$ws.Range("D11").Value = "'0.0912"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -2.96%  "
$ws.Range("D12").Value = "'54.45"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.41%  "
$ws.Range("D13").Value = "'8.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -4.14%  "
$ws.Range("D14").Value = "'0.979"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +8.44%  "
$ws.Range("E15").Value = "  -2.65%  "
$ws.Range("D16").Value = "'14.94"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.22%  "
$ws.Range("D17").Value = "'2.550.94"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.09%  "
$ws.Range("D18").Value = "'2.229.80"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.67%  "
$ws.Range("D19").Value = "'42.219.92"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.95%  "
$ws.Range("D20").Value = "'7.38"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +7.74%  "
$ws.Range("E21").Value = "  -3.99%  "
$ws.Range("D22").Value = "'72.37"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.65%  "
$ws.Range("D23").Value = "'3.47"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +21.49%  "
$ws.Range("E24").Value = "  -3.84%  "
$ws.Range("D25").Value = "'228.10"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.81%  "
$ws.Range("D26").Value = "'9.06"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -3.92%  "
$ws.Range("E27").Value = "  -1.73%  "
$ws.Range("E28").Value = "  -2.11%  "
$ws.Range("E29").Value = "  -0.93%  "
$ws.Range("D30").Value = "'38.19"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -7.81%  "
$ws.Range("D31").Value = "'3.20"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -4.86%  "
$ws.Range("D32").Value = "'173.73"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +1.19%  "
$ws.Range("D33").Value = "'20.95"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.07%  "
$ws.Range("D34").Value = "'0.0896"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -1.66%  "
$ws.Range("E35").Value = "  -1.20%  "
$ws.Range("E36").Value = "  +12.14%  "
$ws.Range("D37").Value = "'4.34"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  +1.45%  "
$ws.Range("E38").Value = "  -2.89%  "
$ws.Range("E39").Value = "  -1.70%  "
$ws.Range("E40").Value = "  -3.68%  "
$ws.Range("E41").Value = "  -4.30%  "
$ws.Range("D42").Value = "'71.53"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.21%  "
$ws.Range("D43").Value = "'0.231"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -1.77%  "
$ws.Range("E44").Value = "  -0.15%  "
$ws.Range("D45").Value = "'12.55"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.08%  "
$ws.Range("E46").Value = "  -4.51%  "
$ws.Range("E47").Value = "  -6.22%  "
$ws.Range("E48").Value = "  +5.17%  "
$ws.Range("D49").Value = "'103.14"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.41%  "
$ws.Range("E50").Value = "  -1.19%  "
$ws.Range("D51").Value = "'1.62"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.65%  "
